$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the old header values in B1/C1 before we shift columns around.
$oldB1 = $ws.Range("B1").Value2
$oldC1 = $ws.Range("C1").Value2

# Insert two new columns at C (existing column C - the most recent week -
# shifts two places to the right, becoming column E).
$ws.Columns("C:D").Insert()

# New header row: B1 becomes the newest week (Jun_17), C1 is the week that
# was just added (Jun_15), D1 keeps the previous newest-week header that
# used to live in B1 (Jun_13), and E1 keeps the old C1 header (Jun_10).
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"
$ws.Range("D1").Value = $oldB1
$ws.Range("E1").Value = $oldC1

# Fill the two freshly-inserted columns with the "UN" placeholder used
# throughout the sheet for rows that have no rating action that week.
For ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}

# Row 19 (Citigroup) got a new Jun_15 rating action; give it the same
# highlight formatting already used elsewhere on this row (column E).
$ws.Range("E19").Copy()
$ws.Range("C19").PasteSpecial(-4122)
$ws.Range("C19").Value = "6/14/2018,Raises Target,Buy -> Buy,$20.00 -> $21.00"

# Row 22 (BidaskClub) also got a new Jun_15 rating action, highlighted the
# same way.
$ws.Range("E19").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value = "6/15/2018,Upgrades,Hold -> Buy,"

# Match the column widths used by the rest of the table for the two new
# columns (same visible width, 8 characters, as the original column C).
$ws.Columns("C").ColumnWidth = 7.17
$ws.Columns("D").ColumnWidth = 7.17

$excel.CutCopyMode = $false
